# Reposition the six plot pictures on slide 3 ("fixed order of plots").
#
# The target positions come from the OOXML diff (EMU). PowerPoint's COM
# object model exposes Shape.Left/.Top in points (1 pt = 12700 EMU) and
# stores them as single-precision floats, so a naive EMU/12700 conversion
# can land one EMU below the intended value after the float32 round-trip.
# The literal point values below were chosen (each is itself an exact
# float32 value) so that Left/Top round-trip to precisely the EMU offsets
# from the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

function Set-ShapePos {
    param($shapes, [string]$name, [double]$xPt, [double]$yPt)
    $sh = $shapes.Item($name)
    $sh.Left = $xPt
    $sh.Top = $yPt
}

# Name         -> target (x EMU, y EMU) == (x pt, y pt)
Set-ShapePos $s.Shapes "Picture 3"  139.875            323.6647644042969   # 1776412, 4110542
Set-ShapePos $s.Shapes "Picture 6"  431.6250915527344  324.8582458496094   # 5481638, 4125699
Set-ShapePos $s.Shapes "Picture 8"  718.8750610351562  324.8583068847656   # 9129713, 4125700
Set-ShapePos $s.Shapes "Picture 22" 113.07854461669922 66.73587036132812  # 1436097, 847545
Set-ShapePos $s.Shapes "Picture 24" 404.8287048339844  66.73587036132812  # 5141324, 847545
Set-ShapePos $s.Shapes "Picture 26" 692.0787353515625  77.64342498779297  # 8789399, 986071
